$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data in rows 2 and 3 (columns B through H); row 2 becomes the
# data that used to be in row 3, and vice versa. Columns A, I, J stay put.
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2
$h2 = $ws.Range("H2").Value2

$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2

$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3
$ws.Range("E2").Value = $e3
$ws.Range("F2").Value = $f3
$ws.Range("G2").Value = $g3
$ws.Range("H2").Value = $h3

$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2
$ws.Range("E3").Value = $e2
$ws.Range("F3").Value = $f2
$ws.Range("G3").Value = $g2
$ws.Range("H3").Value = $h2
